$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 27

# Plain text fields
$ws.Cells.Item($row, 1).Value = "Amirov Akrom Eshali ogli"
$ws.Cells.Item($row, 2).Value = "Hayot faoliyati xavfsizligi"
$ws.Cells.Item($row, 3).Value = "O'zbek tili"
$ws.Cells.Item($row, 4).Value = "Kunduzgi"
$ws.Cells.Item($row, 5).Value = "AB5554414"

# Numeric-looking / date-looking fields: prefix with an apostrophe so
# Excel stores them as text (matching the source inlineStr cells) instead
# of converting them into numbers or date serials.
$ws.Cells.Item($row, 6).Value = "'12345678901234"
$ws.Cells.Item($row, 7).Value = "Xorazm viloyati"
$ws.Cells.Item($row, 8).Value = "Urganch shahri"
$ws.Cells.Item($row, 9).Value = "'+12676860109"
$ws.Cells.Item($row, 10).Value = "'+998945289910"
$ws.Cells.Item($row, 11).Value = "'2025-06-25"
